$d = $word.ActiveDocument
$d.Content.Find.Execute("post-European", $true, $false, $false, $false, $false,
                         $true, 1, $false, "17th century European presence", 2)
